$wb = $excel.ActiveWorkbook

# ALC @@ -1480,22 +1480,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2370.9  # H17: 2269.923 -> 2370.9
$ws.Cells.Item(17, 10).Value = 2370.9  # J17: 2269.923 -> 2370.9
$ws.Cells.Item(17, 12).Value = 7112.700000000001  # L17: 6809.768999999999 -> 7112.700000000001
$ws.Cells.Item(17, 14).Value = -7448.700000000001  # N17: -7145.768999999999 -> -7448.700000000001

# ALC @@ -2610,25 +2610,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3473.2942  # H40: 4040.8696 -> 3473.2942
$ws.Cells.Item(40, 9).Value = 2543.375  # I40: 2535.4285 -> 2543.375
$ws.Cells.Item(40, 10).Value = 4299.8887  # J40: 4699.5 -> 4299.8887
$ws.Cells.Item(40, 11).Value = 2543.375  # K40: 2535.4285 -> 2543.375
$ws.Cells.Item(40, 12).Value = 4299.8887  # L40: 4699.5 -> 4299.8887
$ws.Cells.Item(40, 13).Value = -2368.375  # M40: -2360.4285 -> -2368.375
$ws.Cells.Item(40, 14).Value = -4649.8887  # N40: -5049.5 -> -4649.8887

# ALC @@ -4073,20 +4073,23 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 200006400  # H69: 333339000 -> 200006400
$ws.Cells.Item(69, 9).Value = 7500  # I69: 0 -> 7500
$ws.Cells.Item(69, 11).Value = 22500  # K69: 0 -> 22500
$ws.Cells.Item(69, 13).Value = -21626  # M69: None -> -21626

# ALC @@ -4223,20 +4226,23 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 200006400  # H72: 333339000 -> 200006400
$ws.Cells.Item(72, 9).Value = 7500  # I72: 0 -> 7500
$ws.Cells.Item(72, 11).Value = 67500  # K72: 0 -> 67500
$ws.Cells.Item(72, 13).Value = -63132  # M72: None -> -63132

# ALC @@ -4918,22 +4924,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 5974.25  # H86: 7299.6665 -> 5974.25
$ws.Cells.Item(86, 9).Value = 1998.5  # I86: 1999 -> 1998.5
$ws.Cells.Item(86, 11).Value = 1998.5  # K86: 1999 -> 1998.5
$ws.Cells.Item(86, 13).Value = -875.5  # M86: -876 -> -875.5

# ALC @@ -5071,22 +5077,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 5974.25  # H89: 7299.6665 -> 5974.25
$ws.Cells.Item(89, 9).Value = 1998.5  # I89: 1999 -> 1998.5
$ws.Cells.Item(89, 11).Value = 9992.5  # K89: 9995 -> 9992.5
$ws.Cells.Item(89, 13).Value = -4376.5  # M89: -4379 -> -4376.5

# ALC @@ -5224,25 +5230,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 2609.6667  # H92: 1660.2778 -> 2609.6667
$ws.Cells.Item(92, 9).Value = 3073.8572  # I92: 1603 -> 3073.8572
$ws.Cells.Item(92, 10).Value = 1959.8  # J92: 1774.8334 -> 1959.8
$ws.Cells.Item(92, 11).Value = 3073.8572  # K92: 1603 -> 3073.8572
$ws.Cells.Item(92, 12).Value = 1959.8  # L92: 1774.8334 -> 1959.8
$ws.Cells.Item(92, 13).Value = -1825.8572  # M92: -355 -> -1825.8572
$ws.Cells.Item(92, 14).Value = -4455.8  # N92: -4270.8334 -> -4455.8

# ALC @@ -6179,22 +6185,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 3125.4666  # H111: 3163 -> 3125.4666
$ws.Cells.Item(111, 9).Value = 3085.9167  # I111: 3130.0908 -> 3085.9167
$ws.Cells.Item(111, 11).Value = 9257.750100000001  # K111: 9390.2724 -> 9257.750100000001
$ws.Cells.Item(111, 13).Value = -6190.750100000001  # M111: -6323.2724 -> -6190.750100000001

# ALC @@ -7232,22 +7238,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1657.4054  # H132: 1689.5555 -> 1657.4054
$ws.Cells.Item(132, 9).Value = 1542.4482  # I132: 1579.6786 -> 1542.4482
$ws.Cells.Item(132, 11).Value = 4627.3446  # K132: 4739.0358 -> 4627.3446
$ws.Cells.Item(132, 13).Value = -2097.3446  # M132: -2209.0358 -> -2097.3446

# ALC @@ -7483,25 +7489,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1782.0952  # H137: 1944.9375 -> 1782.0952
$ws.Cells.Item(137, 9).Value = 1823.3889  # I137: 1974.6 -> 1823.3889
$ws.Cells.Item(137, 10).Value = 1534.3334  # J137: 1500 -> 1534.3334
$ws.Cells.Item(137, 11).Value = 5470.1667  # K137: 5923.799999999999 -> 5470.1667
$ws.Cells.Item(137, 12).Value = 4603.0002  # L137: 4500 -> 4603.0002
$ws.Cells.Item(137, 13).Value = -2920.1667  # M137: -3373.799999999999 -> -2920.1667
$ws.Cells.Item(137, 14).Value = -9703.0002  # N137: -9600 -> -9703.0002

# ALC @@ -7685,22 +7691,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 8454.049999999999  # H141: 8199.096 -> 8454.049999999999
$ws.Cells.Item(141, 9).Value = 7672.3335  # I141: 7386.5625 -> 7672.3335
$ws.Cells.Item(141, 11).Value = 23017.0005  # K141: 22159.6875 -> 23017.0005
$ws.Cells.Item(141, 13).Value = -17837.0005  # M141: -16979.6875 -> -17837.0005

# ARM @@ -9325,22 +9331,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7151.7334  # H32: 7356.9653 -> 7151.7334
$ws.Cells.Item(32, 9).Value = 6363.8623  # I32: 6548.2856 -> 6363.8623
$ws.Cells.Item(32, 11).Value = 6363.8623  # K32: 6548.2856 -> 6363.8623
$ws.Cells.Item(32, 13).Value = -6076.8623  # M32: -6261.2856 -> -6076.8623

# ARM @@ -9962,25 +9968,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 9822.571  # H45: 1875.3334 -> 9822.571
$ws.Cells.Item(45, 9).Value = 13270.111  # I45: 1532.35 -> 13270.111
$ws.Cells.Item(45, 10).Value = 3617  # J45: 2855.2856 -> 3617
$ws.Cells.Item(45, 11).Value = 13270.111  # K45: 1532.35 -> 13270.111
$ws.Cells.Item(45, 12).Value = 3617  # L45: 2855.2856 -> 3617
$ws.Cells.Item(45, 13).Value = -12893.111  # M45: -1155.35 -> -12893.111
$ws.Cells.Item(45, 14).Value = -4371  # N45: -3609.2856 -> -4371

# ARM @@ -12351,22 +12357,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(94, 8).Value = 73776.664  # H94: 30315 -> 73776.664
$ws.Cells.Item(94, 9).Value = 100000  # I94: 0 -> 100000
$ws.Cells.Item(94, 10).Value = 60665  # J94: 30315 -> 60665
$ws.Cells.Item(94, 11).Value = 100000  # K94: 0 -> 100000
$ws.Cells.Item(94, 12).Value = 60665  # L94: 30315 -> 60665
$ws.Cells.Item(94, 13).Value = -99099  # M94: None -> -99099
$ws.Cells.Item(94, 14).Value = -62467  # N94: -32117 -> -62467

# BSM @@ -15703,25 +15712,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 6014.6523  # H20: 6253.091 -> 6014.6523
$ws.Cells.Item(20, 9).Value = 8182.1333  # I20: 8184.8 -> 8182.1333
$ws.Cells.Item(20, 10).Value = 1950.625  # J20: 2113.7144 -> 1950.625
$ws.Cells.Item(20, 11).Value = 8182.1333  # K20: 8184.8 -> 8182.1333
$ws.Cells.Item(20, 12).Value = 1950.625  # L20: 2113.7144 -> 1950.625
$ws.Cells.Item(20, 13).Value = -7935.1333  # M20: -7937.8 -> -7935.1333
$ws.Cells.Item(20, 14).Value = -2444.625  # N20: -2607.7144 -> -2444.625

# BSM @@ -15801,25 +15810,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 11355.444  # H22: 14547.143 -> 11355.444
$ws.Cells.Item(22, 9).Value = 25275  # I22: 20266 -> 25275
$ws.Cells.Item(22, 10).Value = 219.8  # J22: 250 -> 219.8
$ws.Cells.Item(22, 11).Value = 25275  # K22: 20266 -> 25275
$ws.Cells.Item(22, 12).Value = 219.8  # L22: 250 -> 219.8
$ws.Cells.Item(22, 13).Value = -25102  # M22: -20093 -> -25102
$ws.Cells.Item(22, 14).Value = -565.8  # N22: -596 -> -565.8

# BSM @@ -19305,25 +19314,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1859.1111  # H94: 2539.5 -> 1859.1111
$ws.Cells.Item(94, 9).Value = 2364.5  # I94: 3296.75 -> 2364.5
$ws.Cells.Item(94, 10).Value = 848.3333  # J94: 1025 -> 848.3333
$ws.Cells.Item(94, 11).Value = 2364.5  # K94: 3296.75 -> 2364.5
$ws.Cells.Item(94, 12).Value = 848.3333  # L94: 1025 -> 848.3333
$ws.Cells.Item(94, 13).Value = -1913.5  # M94: -2845.75 -> -1913.5
$ws.Cells.Item(94, 14).Value = -1750.3333  # N94: -1927 -> -1750.3333

# BSM @@ -21262,22 +21271,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2545.3235  # H134: 2595.2424 -> 2545.3235
$ws.Cells.Item(134, 9).Value = 2417.0908  # I134: 2489.4285 -> 2417.0908
$ws.Cells.Item(134, 11).Value = 7251.2724  # K134: 7468.2855 -> 7251.2724
$ws.Cells.Item(134, 13).Value = -4716.2724  # M134: -4933.2855 -> -4716.2724

# CRP @@ -22455,25 +22464,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1951.625  # H16: 2166.6667 -> 1951.625
$ws.Cells.Item(16, 9).Value = 1760  # I16: 1750 -> 1760
$ws.Cells.Item(16, 10).Value = 2271  # J16: 3000 -> 2271
$ws.Cells.Item(16, 11).Value = 1760  # K16: 1750 -> 1760
$ws.Cells.Item(16, 12).Value = 2271  # L16: 3000 -> 2271
$ws.Cells.Item(16, 13).Value = -1473  # M16: -1463 -> -1473
$ws.Cells.Item(16, 14).Value = -2845  # N16: -3574 -> -2845

# CRP @@ -23187,22 +23196,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3423.4707  # H31: 3520.7878 -> 3423.4707
$ws.Cells.Item(31, 9).Value = 1757.909  # I31: 1831.5238 -> 1757.909
$ws.Cells.Item(31, 11).Value = 1757.909  # K31: 1831.5238 -> 1757.909
$ws.Cells.Item(31, 13).Value = -1462.909  # M31: -1536.5238 -> -1462.909

# CRP @@ -23340,22 +23349,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3423.4707  # H34: 3520.7878 -> 3423.4707
$ws.Cells.Item(34, 9).Value = 1757.909  # I34: 1831.5238 -> 1757.909
$ws.Cells.Item(34, 11).Value = 1757.909  # K34: 1831.5238 -> 1757.909
$ws.Cells.Item(34, 13).Value = -1555.909  # M34: -1629.5238 -> -1555.909

# CRP @@ -27247,25 +27256,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1951.625  # H113: 2166.6667 -> 1951.625
$ws.Cells.Item(113, 9).Value = 1760  # I113: 1750 -> 1760
$ws.Cells.Item(113, 10).Value = 2271  # J113: 3000 -> 2271
$ws.Cells.Item(113, 11).Value = 1760  # K113: 1750 -> 1760
$ws.Cells.Item(113, 12).Value = 2271  # L113: 3000 -> 2271
$ws.Cells.Item(113, 13).Value = 410  # M113: 420 -> 410
$ws.Cells.Item(113, 14).Value = -6611  # N113: -7340 -> -6611

# CUL @@ -35275,25 +35284,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 8100  # H131: 3453.5386 -> 8100
$ws.Cells.Item(131, 9).Value = 1500  # I131: 1333.3334 -> 1500
$ws.Cells.Item(131, 10).Value = 12500  # J131: 4089.6 -> 12500
$ws.Cells.Item(131, 11).Value = 4500  # K131: 4000.0002 -> 4500
$ws.Cells.Item(131, 12).Value = 37500  # L131: 12268.8 -> 37500
$ws.Cells.Item(131, 13).Value = 540  # M131: 1039.9998 -> 540
$ws.Cells.Item(131, 14).Value = -47580  # N131: -22348.8 -> -47580

# CUL @@ -35688,22 +35697,22 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 4350.7085  # H139: 4512.478 -> 4350.7085
$ws.Cells.Item(139, 9).Value = 3801.3076  # I139: 4065.5833 -> 3801.3076
$ws.Cells.Item(139, 11).Value = 11403.9228  # K139: 12196.7499 -> 11403.9228
$ws.Cells.Item(139, 13).Value = -6263.9228  # M139: -7056.749899999999 -> -6263.9228

# GSM @@ -36783,26 +36792,23 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 4999  # H19: 2666.3333 -> 4999
$ws.Cells.Item(19, 10).Value = 0  # J19: 1500 -> 0
$ws.Cells.Item(19, 12).Value = 0  # L19: 1500 -> 0
$ws.Cells.Item(19, 14).ClearContents()  # N19: -2076 -> (removed)

# GSM @@ -38268,22 +38274,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 26249.75  # H49: 27500.5 -> 26249.75
$ws.Cells.Item(49, 10).Value = 26249.75  # J49: 27500.5 -> 26249.75
$ws.Cells.Item(49, 12).Value = 26249.75  # L49: 27500.5 -> 26249.75
$ws.Cells.Item(49, 14).Value = -26617.75  # N49: -27868.5 -> -26617.75

# GSM @@ -39291,22 +39297,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 11965.947  # H70: 13031.588 -> 11965.947
$ws.Cells.Item(70, 9).Value = 12523.818  # I70: 14660.667 -> 12523.818
$ws.Cells.Item(70, 11).Value = 12523.818  # K70: 14660.667 -> 12523.818
$ws.Cells.Item(70, 13).Value = -12253.818  # M70: -14390.667 -> -12253.818

# GSM @@ -39441,22 +39447,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 11965.947  # H73: 13031.588 -> 11965.947
$ws.Cells.Item(73, 9).Value = 12523.818  # I73: 14660.667 -> 12523.818
$ws.Cells.Item(73, 11).Value = 12523.818  # K73: 14660.667 -> 12523.818
$ws.Cells.Item(73, 13).Value = -11587.818  # M73: -13724.667 -> -11587.818

# GSM @@ -40427,25 +40433,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 22428.143  # H93: 21874.5 -> 22428.143
$ws.Cells.Item(93, 10).Value = 23332.834  # J93: 22570.857 -> 23332.834
$ws.Cells.Item(93, 12).Value = 23332.834  # L93: 22570.857 -> 23332.834
$ws.Cells.Item(93, 14).Value = -27076.834  # N93: -26314.857 -> -27076.834

# GSM @@ -40620,25 +40626,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1073.5454  # H97: 1132.8 -> 1073.5454
$ws.Cells.Item(97, 9).Value = 1152.1578  # I97: 1205.1111 -> 1152.1578
$ws.Cells.Item(97, 10).Value = 575.6667  # J97: 482 -> 575.6667
$ws.Cells.Item(97, 11).Value = 1152.1578  # K97: 1205.1111 -> 1152.1578
$ws.Cells.Item(97, 12).Value = 575.6667  # L97: 482 -> 575.6667
$ws.Cells.Item(97, 13).Value = -656.1578  # M97: -709.1111000000001 -> -656.1578
$ws.Cells.Item(97, 14).Value = -1567.6667  # N97: -1474 -> -1567.6667

# GSM @@ -41113,25 +41119,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 1414.3914  # H107: 1523.6666 -> 1414.3914
$ws.Cells.Item(107, 9).Value = 366.18182  # I107: 425.5 -> 366.18182
$ws.Cells.Item(107, 10).Value = 2375.25  # J107: 2199.4614 -> 2375.25
$ws.Cells.Item(107, 11).Value = 366.18182  # K107: 425.5 -> 366.18182
$ws.Cells.Item(107, 12).Value = 2375.25  # L107: 2199.4614 -> 2375.25
$ws.Cells.Item(107, 13).Value = 1553.81818  # M107: 1494.5 -> 1553.81818
$ws.Cells.Item(107, 14).Value = -6215.25  # N107: -6039.4614 -> -6215.25

# LTW @@ -43902,25 +43908,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2178.5625  # H22: 2329.5833 -> 2178.5625
$ws.Cells.Item(22, 9).Value = 1650.1666  # I22: 1667 -> 1650.1666
$ws.Cells.Item(22, 10).Value = 2495.6  # J22: 2550.4443 -> 2495.6
$ws.Cells.Item(22, 11).Value = 1650.1666  # K22: 1667 -> 1650.1666
$ws.Cells.Item(22, 12).Value = 2495.6  # L22: 2550.4443 -> 2495.6
$ws.Cells.Item(22, 13).Value = -1355.1666  # M22: -1372 -> -1355.1666
$ws.Cells.Item(22, 14).Value = -3085.6  # N22: -3140.4443 -> -3085.6

# LTW @@ -44144,25 +44150,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 2178.5625  # H27: 2329.5833 -> 2178.5625
$ws.Cells.Item(27, 9).Value = 1650.1666  # I27: 1667 -> 1650.1666
$ws.Cells.Item(27, 10).Value = 2495.6  # J27: 2550.4443 -> 2495.6
$ws.Cells.Item(27, 11).Value = 1650.1666  # K27: 1667 -> 1650.1666
$ws.Cells.Item(27, 12).Value = 2495.6  # L27: 2550.4443 -> 2495.6
$ws.Cells.Item(27, 13).Value = -1543.1666  # M27: -1560 -> -1543.1666
$ws.Cells.Item(27, 14).Value = -2709.6  # N27: -2764.4443 -> -2709.6

# LTW @@ -44888,20 +44894,23 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 17833.334  # H42: 21000 -> 17833.334
$ws.Cells.Item(42, 9).Value = 2000  # I42: 0 -> 2000
$ws.Cells.Item(42, 11).Value = 2000  # K42: 0 -> 2000
$ws.Cells.Item(42, 13).Value = -1437  # M42: None -> -1437

# LTW @@ -45081,25 +45090,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1720.9375  # H46: 1649.7059 -> 1720.9375
$ws.Cells.Item(46, 10).Value = 2127.9167  # J46: 2128.75 -> 2127.9167
$ws.Cells.Item(46, 12).Value = 2127.9167  # L46: 2128.75 -> 2127.9167
$ws.Cells.Item(46, 14).Value = -2503.9167  # N46: -2504.75 -> -2503.9167

# LTW @@ -45225,20 +45234,23 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(49, 8).Value = 17833.334  # H49: 21000 -> 17833.334
$ws.Cells.Item(49, 9).Value = 2000  # I49: 0 -> 2000
$ws.Cells.Item(49, 11).Value = 2000  # K49: 0 -> 2000
$ws.Cells.Item(49, 13).Value = -1853  # M49: None -> -1853

# LTW @@ -48814,22 +48826,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3521.8125  # H122: 3484.647 -> 3521.8125
$ws.Cells.Item(122, 9).Value = 3353.182  # I122: 3314.5833 -> 3353.182
$ws.Cells.Item(122, 11).Value = 10059.546  # K122: 9943.749899999999 -> 10059.546
$ws.Cells.Item(122, 13).Value = -7609.545999999998  # M122: -7493.749899999999 -> -7609.545999999998

# WVR @@ -52011,22 +52023,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 48722  # H46: 49970.75 -> 48722
$ws.Cells.Item(46, 10).Value = 48722  # J46: 49970.75 -> 48722
$ws.Cells.Item(46, 12).Value = 48722  # L46: 49970.75 -> 48722
$ws.Cells.Item(46, 14).Value = -49184  # N46: -50432.75 -> -49184

# WVR @@ -53738,22 +53750,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2413.875  # H81: 2109 -> 2413.875
$ws.Cells.Item(81, 9).Value = 1684.5  # I81: 1419.5 -> 1684.5
$ws.Cells.Item(81, 11).Value = 3369  # K81: 2839 -> 3369
$ws.Cells.Item(81, 13).Value = -2308  # M81: -1778 -> -2308

# WVR @@ -53885,22 +53897,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 2413.875  # H84: 2109 -> 2413.875
$ws.Cells.Item(84, 9).Value = 1684.5  # I84: 1419.5 -> 1684.5
$ws.Cells.Item(84, 11).Value = 16845  # K84: 14195 -> 16845
$ws.Cells.Item(84, 13).Value = -11541  # M84: -8891 -> -11541

# WVR @@ -55015,22 +55027,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 628.41174  # H107: 621.2778 -> 628.41174
$ws.Cells.Item(107, 9).Value = 527.5  # I107: 525.38464 -> 527.5
$ws.Cells.Item(107, 11).Value = 1582.5  # K107: 1576.15392 -> 1582.5
$ws.Cells.Item(107, 13).Value = 337.5  # M107: 343.84608 -> 337.5

# WVR @@ -55744,22 +55756,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3416.1365  # H122: 3488.6047 -> 3416.1365
$ws.Cells.Item(122, 9).Value = 2033.8  # I122: 2157.6428 -> 2033.8
$ws.Cells.Item(122, 11).Value = 6101.4  # K122: 6472.928400000001 -> 6101.4
$ws.Cells.Item(122, 13).Value = -3651.4  # M122: -4022.928400000001 -> -3651.4

# WVR @@ -56228,22 +56240,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2946.1482  # H132: 3021.1924 -> 2946.1482
$ws.Cells.Item(132, 9).Value = 2454.0952  # I132: 2527.05 -> 2454.0952
$ws.Cells.Item(132, 11).Value = 7362.285600000001  # K132: 7581.150000000001 -> 7362.285600000001
$ws.Cells.Item(132, 13).Value = -4832.285600000001  # M132: -5051.150000000001 -> -4832.285600000001

# WVR @@ -56329,22 +56341,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(134, 8).Value = 48722  # H134: 49970.75 -> 48722
$ws.Cells.Item(134, 10).Value = 48722  # J134: 49970.75 -> 48722
$ws.Cells.Item(134, 12).Value = 146166  # L134: 149912.25 -> 146166
$ws.Cells.Item(134, 14).Value = -151236  # N134: -154982.25 -> -151236
